$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(312049950, Molham  Peretz: 7,2)"
$ws.Range("B1").Value = "(308073899, Anan  Kirshenbaum: 5,-1)"
$ws.Range("C1").Value = "(318869187, Soaad  Leibovich: 4,-4)"
$ws.Range("D1").Value = "(205898513, Asaf  Braymok: 3,-4)"
$ws.Range("E1").Value = "(318294931, Shalev  Afanasenko: -7,-4)"
$ws.Range("F1").Value = "(318428158, Tal  Asulin: -6,2)"
$ws.Range("G1").Value = "(316028364, Sami  Castro: -1,6)"

$ws.Range("A3").Value = "cost: 253.644346252802"
$ws.Range("A4").Value = "time: 45.72886925056039"
